$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 324, shifting existing rows 324-347 down to 325-348
$ws.Rows(324).Insert()

# Populate the newly inserted row 324 with the new record
$ws.Cells.Item(324, 1).Value = 3
$ws.Cells.Item(324, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(324, 3).Value = "Coquimbo"
$ws.Cells.Item(324, 4).Value = 44714
$ws.Cells.Item(324, 5).Value = 5
$ws.Cells.Item(324, 6).Value = 100112012
$ws.Cells.Item(324, 7).Value = "Espinaca"
$ws.Cells.Item(324, 8).Value = "Sin especificar"
$ws.Cells.Item(324, 9).Value = "Primera"
$ws.Cells.Item(324, 10).Value = 150
$ws.Cells.Item(324, 11).Value = 4000
$ws.Cells.Item(324, 12).Value = 4000
$ws.Cells.Item(324, 13).Value = 4000
$ws.Cells.Item(324, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(324, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(324, 16).Value = 1333
$ws.Cells.Item(324, 17).Value = 3
$ws.Cells.Item(324, 18).Value = "Hortaliza"
